# Update crypto price/volume data (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "23.923.17"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  -0.66%  "

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "1.650.85"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  +0.84%  "

$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = "1.000"
$dCell.Style = "Normal"
$ws.Range("E4").Value = "  +0.23%  "

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "308.97"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  -0.84%  "

$ws.Range("E6").Value = "  +0.28%  "

$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = "0.3896"
$dCell.Style = "Normal"
$ws.Range("E7").Value = "  -1.21%  "

$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "0.3833"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  -1.07%  "

$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = "51.84"
$dCell.Style = "Normal"
$ws.Range("E9").Value = "  -0.63%  "

$ws.Range("E10").Value = "  -2.15%  "

$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "1.002"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  +0.34%  "

$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = "0.08428"
$dCell.Style = "Normal"
$ws.Range("E12").Value = "  -1.12%  "

$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "23.88"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  -0.92%  "

$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = "7.083"
$dCell.Style = "Normal"
$ws.Range("E14").Value = "  -0.75%  "

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "7.955"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  +3.52%  "

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "0.00001317"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  +1.21%  "

$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "1.646.58"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  +2.00%  "

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "94.76"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  -0.01%  "

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "0.06969"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  +0.76%  "

$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "19.72"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  -2.41%  "

$ws.Range("E21").Value = "  +0.58%  "

$ws.Range("E22").Value = "  +0.44%  "

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "13.73"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  +1.52%  "

$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "23.923.33"
$dCell.Style = "Normal"
$ws.Range("E24").Value = "  -0.61%  "

$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = "2.454"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  -0.11%  "

$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = "2.956"
$dCell.Style = "Normal"
$ws.Range("E26").Value = "  +2.28%  "

$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = "22.10"
$dCell.Style = "Normal"
$ws.Range("E27").Value = "  -1.19%  "

$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "151.47"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  -3.67%  "

$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = "5.403"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  +0.88%  "

$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "139.00"
$dCell.Style = "Normal"
$ws.Range("E30").Value = "  -1.52%  "

$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "7.877"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  -1.98%  "

$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "2.518"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  -0.47%  "

$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = "1.829.12"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  +1.74%  "

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "1.046"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  +3.44%  "

$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "0.08032"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  -1.81%  "

$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "0.02965"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  +1.16%  "

$ws.Range("E37").Value = "  +4.70%  "

$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "6.673"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  -0.40%  "

$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = "0.2682"
$dCell.Style = "Normal"
$ws.Range("E39").Value = "  -0.42%  "

$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "0.09101"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  -1.17%  "

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "0.7618"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("E42").Value = "  -2.66%  "

$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "1.425"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  -0.55%  "

$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = "16.31"
$dCell.Style = "Normal"
$ws.Range("E44").Value = "  +0.27%  "

$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "0.7007"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  +0.52%  "

$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = "2.468"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  -0.83%  "

$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "4.077"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  -0.73%  "

$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "0.08292"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  -0.73%  "

$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "134.53"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  -1.49%  "

$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "1.214"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  -1.74%  "
